# Change the "NIDN : " label (in front of the {nidn_nuptk} placeholder)
# to "NUPTK" + " : ", split across two runs, matching the target diff
# while leaving the following "{nidn_nuptk}" run untouched.

$d = $word.ActiveDocument

# Find the paragraph that holds the literal "NIDN : {nidn_nuptk}" text --
# there is another unrelated "NIDN/NUPTK" occurrence in the footer, so we
# match on the unique, longer substring to make sure we grab the right
# paragraph.
$target = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $para = $paras.Item($i)
    if ($para.Range.Text -like "NIDN : {nidn_nuptk}*") {
        $target = $para
    }
}

$r = $target.Range
$start = $r.Start

# Replace "NIDN : " (7 characters) with "NUPTK : ".
$editRange = $d.Range($start, $start + 7)
$editRange.Text = "NUPTK : "

# The text replacement above collapses the whole paragraph into a single
# run. Force Word to re-split it into separate runs at the two boundaries
# we need ("NUPTK" | " : " | "{nidn_nuptk}") by toggling a character
# attribute on/off across just the leading sub-ranges -- toggling back to
# the original value keeps formatting identical while leaving a run
# boundary behind.
$boundary1 = $d.Range($start, $start + 5)
$boundary1.Bold = 1
$boundary1.Bold = 0

$boundary2 = $d.Range($start, $start + 8)
$boundary2.Bold = 1
$boundary2.Bold = 0
